$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New lawn-sign rows appended to the bottom of the tracking sheet
# (Guildwood Way / Starwood Dr / Ceremonial Dr addresses), matching the
# "checking in AP08 file and small lawn sign" commit.

$rows = @(
    @{ Row = 45; Num = 5256; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 46; Num = 5076; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 47; Num = 5088; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 48; Num = 5104; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 49; Num = 5160; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 50; Num = 5194; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 51; Num = 5240; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 52; Num = 5036; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 53; Num = 5096; Street = "Guildwood Way";  Installed = "Y" },
    @{ Row = 54; Num = 5447; Street = "Starwood Dr";    Installed = $null },
    @{ Row = 55; Num = 720;  Street = "Ceremonial Dr";  Installed = $null }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Num

    $c = $ws.Cells.Item($r.Row, 3)
    $c.NumberFormat = "@"
    $c.Value = $r.Street

    if ($r.Installed) {
        $ws.Cells.Item($r.Row, 6).Value = $r.Installed
    }
}

# Match the author's final selection/scroll state after entering the rows
[void]$ws.Cells.Item(55, 2).Select()

Write-Host "Appended lawn sign rows 45-55"
